# Apply the match-row rotation described in the commit diff.
# Each block of rows keeps its sequential index (column A) fixed,
# while the match data (columns B, E:AD) rotates by one position
# (last row's data moves to the first row of the block).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 124
$ws.Cells.Item(124, 2).Value = 6936857
$ws.Cells.Item(124, 5).Value = 'AEK Athens'
$ws.Cells.Item(124, 6).Value = 'Panathinaikos'
$ws.Cells.Item(124, 7).Value = 2
$ws.Cells.Item(124, 8).Value = 2
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = 1
$ws.Cells.Item(124, 11).Value = 'D'
$ws.Cells.Item(124, 12).Value = 1.909
$ws.Cells.Item(124, 13).Value = 3.5
$ws.Cells.Item(124, 14).Value = 4.2
$ws.Cells.Item(124, 15).Value = 2.15
$ws.Cells.Item(124, 16).Value = 3.2
$ws.Cells.Item(124, 17).Value = 3.5
$ws.Cells.Item(124, 18).Value = -0.25
$ws.Cells.Item(124, 19).Value = 1.85
$ws.Cells.Item(124, 20).Value = 2
$ws.Cells.Item(124, 21).Value = 2
$ws.Cells.Item(124, 22).Value = 1.8
$ws.Cells.Item(124, 23).Value = 2.05
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 25).Value = 2.2
$ws.Cells.Item(124, 26).Value = -1
$ws.Cells.Item(124, 27).Value = -0.5
$ws.Cells.Item(124, 28).Value = 0.5
$ws.Cells.Item(124, 29).Value = 0.8
$ws.Cells.Item(124, 30).Value = -1

# Row 125
$ws.Cells.Item(125, 2).Value = 6937238
$ws.Cells.Item(125, 5).Value = 'PAOK Salonika'
$ws.Cells.Item(125, 6).Value = 'Giannina'
$ws.Cells.Item(125, 7).Value = 4
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 'H'
$ws.Cells.Item(125, 12).Value = 1.111
$ws.Cells.Item(125, 13).Value = 9
$ws.Cells.Item(125, 14).Value = 23
$ws.Cells.Item(125, 15).Value = 1.25
$ws.Cells.Item(125, 16).Value = 6
$ws.Cells.Item(125, 17).Value = 9
$ws.Cells.Item(125, 18).Value = -1.75
$ws.Cells.Item(125, 19).Value = 2.025
$ws.Cells.Item(125, 20).Value = 1.825
$ws.Cells.Item(125, 21).Value = 2.75
$ws.Cells.Item(125, 22).Value = 1.8
$ws.Cells.Item(125, 23).Value = 2.05
$ws.Cells.Item(125, 24).Value = 0.25
$ws.Cells.Item(125, 25).Value = -1
$ws.Cells.Item(125, 26).Value = -1
$ws.Cells.Item(125, 27).Value = 1.025
$ws.Cells.Item(125, 28).Value = -1
$ws.Cells.Item(125, 29).Value = 0.8
$ws.Cells.Item(125, 30).Value = -1

# Row 175
$ws.Cells.Item(175, 2).Value = 6935700
$ws.Cells.Item(175, 5).Value = 'Panserraikos'
$ws.Cells.Item(175, 6).Value = 'Asteras Tripolis'
$ws.Cells.Item(175, 7).Value = 2
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(175, 9).Value = 1
$ws.Cells.Item(175, 10).Value = 1
$ws.Cells.Item(175, 11).Value = 'H'
$ws.Cells.Item(175, 12).Value = 2.6
$ws.Cells.Item(175, 13).Value = 3.2
$ws.Cells.Item(175, 14).Value = 2.875
$ws.Cells.Item(175, 15).Value = 2.25
$ws.Cells.Item(175, 16).Value = 3.3
$ws.Cells.Item(175, 17).Value = 3.3
$ws.Cells.Item(175, 18).Value = -0.25
$ws.Cells.Item(175, 19).Value = 1.925
$ws.Cells.Item(175, 20).Value = 1.925
$ws.Cells.Item(175, 21).Value = 2.25
$ws.Cells.Item(175, 22).Value = 2
$ws.Cells.Item(175, 23).Value = 1.85
$ws.Cells.Item(175, 24).Value = 1.25
$ws.Cells.Item(175, 25).Value = -1
$ws.Cells.Item(175, 26).Value = -1
$ws.Cells.Item(175, 27).Value = 0.925
$ws.Cells.Item(175, 28).Value = -1
$ws.Cells.Item(175, 29).Value = 1
$ws.Cells.Item(175, 30).Value = -1

# Row 176
$ws.Cells.Item(176, 2).Value = 6935701
$ws.Cells.Item(176, 5).Value = 'Kifisias FC'
$ws.Cells.Item(176, 6).Value = 'Panetolikos'
$ws.Cells.Item(176, 7).Value = 2
$ws.Cells.Item(176, 8).Value = 2
$ws.Cells.Item(176, 9).Value = 1
$ws.Cells.Item(176, 10).Value = 0
$ws.Cells.Item(176, 11).Value = 'D'
$ws.Cells.Item(176, 12).Value = 2.45
$ws.Cells.Item(176, 13).Value = 3.25
$ws.Cells.Item(176, 14).Value = 3
$ws.Cells.Item(176, 15).Value = 2.05
$ws.Cells.Item(176, 16).Value = 3.3
$ws.Cells.Item(176, 17).Value = 3.8
$ws.Cells.Item(176, 18).Value = -0.5
$ws.Cells.Item(176, 19).Value = 2.05
$ws.Cells.Item(176, 20).Value = 1.8
$ws.Cells.Item(176, 21).Value = 2.25
$ws.Cells.Item(176, 22).Value = 1.8
$ws.Cells.Item(176, 23).Value = 2.05
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 2.3
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = -1
$ws.Cells.Item(176, 28).Value = 0.8
$ws.Cells.Item(176, 29).Value = 0.8
$ws.Cells.Item(176, 30).Value = -1

# Row 177
$ws.Cells.Item(177, 2).Value = 6936863
$ws.Cells.Item(177, 5).Value = 'OFI Crete'
$ws.Cells.Item(177, 6).Value = 'Panathinaikos'
$ws.Cells.Item(177, 7).Value = 2
$ws.Cells.Item(177, 8).Value = 2
$ws.Cells.Item(177, 9).Value = 1
$ws.Cells.Item(177, 10).Value = 1
$ws.Cells.Item(177, 11).Value = 'D'
$ws.Cells.Item(177, 12).Value = 8
$ws.Cells.Item(177, 13).Value = 4.75
$ws.Cells.Item(177, 14).Value = 1.4
$ws.Cells.Item(177, 15).Value = 5.5
$ws.Cells.Item(177, 16).Value = 4.75
$ws.Cells.Item(177, 17).Value = 1.55
$ws.Cells.Item(177, 18).Value = 1
$ws.Cells.Item(177, 19).Value = 1.95
$ws.Cells.Item(177, 20).Value = 1.9
$ws.Cells.Item(177, 21).Value = 2.5
$ws.Cells.Item(177, 22).Value = 1.85
$ws.Cells.Item(177, 23).Value = 2
$ws.Cells.Item(177, 24).Value = -1
$ws.Cells.Item(177, 25).Value = 3.75
$ws.Cells.Item(177, 26).Value = -1
$ws.Cells.Item(177, 27).Value = 0.95
$ws.Cells.Item(177, 28).Value = -1
$ws.Cells.Item(177, 29).Value = 0.8500000000000001
$ws.Cells.Item(177, 30).Value = -1

# Row 178
$ws.Cells.Item(178, 2).Value = 6937269
$ws.Cells.Item(178, 5).Value = 'Aris Salonika'
$ws.Cells.Item(178, 6).Value = 'AEK Athens'
$ws.Cells.Item(178, 7).Value = 3
$ws.Cells.Item(178, 8).Value = 3
$ws.Cells.Item(178, 9).Value = 1
$ws.Cells.Item(178, 10).Value = 1
$ws.Cells.Item(178, 11).Value = 'D'
$ws.Cells.Item(178, 12).Value = 4.75
$ws.Cells.Item(178, 13).Value = 3.75
$ws.Cells.Item(178, 14).Value = 1.75
$ws.Cells.Item(178, 15).Value = 6.5
$ws.Cells.Item(178, 16).Value = 4.2
$ws.Cells.Item(178, 17).Value = 1.5
$ws.Cells.Item(178, 18).Value = 1
$ws.Cells.Item(178, 19).Value = 2.05
$ws.Cells.Item(178, 20).Value = 1.8
$ws.Cells.Item(178, 21).Value = 2.5
$ws.Cells.Item(178, 22).Value = 1.975
$ws.Cells.Item(178, 23).Value = 1.875
$ws.Cells.Item(178, 24).Value = -1
$ws.Cells.Item(178, 25).Value = 3.2
$ws.Cells.Item(178, 26).Value = -1
$ws.Cells.Item(178, 27).Value = 1.05
$ws.Cells.Item(178, 28).Value = -1
$ws.Cells.Item(178, 29).Value = 0.9750000000000001
$ws.Cells.Item(178, 30).Value = -1

# Row 179
$ws.Cells.Item(179, 2).Value = 6937270
$ws.Cells.Item(179, 5).Value = 'Olympiakos'
$ws.Cells.Item(179, 6).Value = 'Volos NFC'
$ws.Cells.Item(179, 7).Value = 3
$ws.Cells.Item(179, 8).Value = 0
$ws.Cells.Item(179, 9).Value = 2
$ws.Cells.Item(179, 10).Value = 0
$ws.Cells.Item(179, 11).Value = 'H'
$ws.Cells.Item(179, 12).Value = 1.125
$ws.Cells.Item(179, 13).Value = 9
$ws.Cells.Item(179, 14).Value = 19
$ws.Cells.Item(179, 15).Value = 1.111
$ws.Cells.Item(179, 16).Value = 9
$ws.Cells.Item(179, 17).Value = 21
$ws.Cells.Item(179, 18).Value = -2.25
$ws.Cells.Item(179, 19).Value = 1.875
$ws.Cells.Item(179, 20).Value = 1.975
$ws.Cells.Item(179, 21).Value = 3.25
$ws.Cells.Item(179, 22).Value = 2
$ws.Cells.Item(179, 23).Value = 1.85
$ws.Cells.Item(179, 24).Value = 0.111
$ws.Cells.Item(179, 25).Value = -1
$ws.Cells.Item(179, 26).Value = -1
$ws.Cells.Item(179, 27).Value = 0.875
$ws.Cells.Item(179, 28).Value = -1
$ws.Cells.Item(179, 29).Value = -0.5
$ws.Cells.Item(179, 30).Value = 0.425

# Row 180
$ws.Cells.Item(180, 2).Value = 6937271
$ws.Cells.Item(180, 5).Value = 'Giannina'
$ws.Cells.Item(180, 6).Value = 'Atromitos Athinon'
$ws.Cells.Item(180, 7).Value = 1
$ws.Cells.Item(180, 8).Value = 1
$ws.Cells.Item(180, 9).Value = 1
$ws.Cells.Item(180, 10).Value = 0
$ws.Cells.Item(180, 11).Value = 'D'
$ws.Cells.Item(180, 12).Value = 2.45
$ws.Cells.Item(180, 13).Value = 3.1
$ws.Cells.Item(180, 14).Value = 3.1
$ws.Cells.Item(180, 15).Value = 2
$ws.Cells.Item(180, 16).Value = 3.3
$ws.Cells.Item(180, 17).Value = 4
$ws.Cells.Item(180, 18).Value = -0.5
$ws.Cells.Item(180, 19).Value = 2.025
$ws.Cells.Item(180, 20).Value = 1.825
$ws.Cells.Item(180, 21).Value = 2.25
$ws.Cells.Item(180, 22).Value = 1.85
$ws.Cells.Item(180, 23).Value = 2
$ws.Cells.Item(180, 24).Value = -1
$ws.Cells.Item(180, 25).Value = 2.3
$ws.Cells.Item(180, 26).Value = -1
$ws.Cells.Item(180, 27).Value = -1
$ws.Cells.Item(180, 28).Value = 0.825
$ws.Cells.Item(180, 29).Value = -0.5
$ws.Cells.Item(180, 30).Value = 0.5

# Row 194
$ws.Cells.Item(194, 2).Value = 7920471
$ws.Cells.Item(194, 5).Value = 'Aris Salonika'
$ws.Cells.Item(194, 6).Value = 'Lamia'
$ws.Cells.Item(194, 7).Value = 3
$ws.Cells.Item(194, 8).Value = 1
$ws.Cells.Item(194, 9).Value = 0
$ws.Cells.Item(194, 10).Value = 0
$ws.Cells.Item(194, 11).Value = 'H'
$ws.Cells.Item(194, 12).Value = 1.571
$ws.Cells.Item(194, 13).Value = 4
$ws.Cells.Item(194, 14).Value = 6
$ws.Cells.Item(194, 15).Value = 1.444
$ws.Cells.Item(194, 16).Value = 4.5
$ws.Cells.Item(194, 17).Value = 8.5
$ws.Cells.Item(194, 18).Value = -1.25
$ws.Cells.Item(194, 19).Value = 1.925
$ws.Cells.Item(194, 20).Value = 1.925
$ws.Cells.Item(194, 21).Value = 2.75
$ws.Cells.Item(194, 22).Value = 2.025
$ws.Cells.Item(194, 23).Value = 1.825
$ws.Cells.Item(194, 24).Value = 0.444
$ws.Cells.Item(194, 25).Value = -1
$ws.Cells.Item(194, 26).Value = -1
$ws.Cells.Item(194, 27).Value = 0.925
$ws.Cells.Item(194, 28).Value = -1
$ws.Cells.Item(194, 29).Value = 1.025
$ws.Cells.Item(194, 30).Value = -1

# Row 195
$ws.Cells.Item(195, 2).Value = 7920470
$ws.Cells.Item(195, 5).Value = 'AEK Athens'
$ws.Cells.Item(195, 6).Value = 'Olympiakos'
$ws.Cells.Item(195, 7).Value = 1
$ws.Cells.Item(195, 8).Value = 0
$ws.Cells.Item(195, 9).Value = 0
$ws.Cells.Item(195, 10).Value = 0
$ws.Cells.Item(195, 11).Value = 'H'
$ws.Cells.Item(195, 12).Value = 1.909
$ws.Cells.Item(195, 13).Value = 3.4
$ws.Cells.Item(195, 14).Value = 4.2
$ws.Cells.Item(195, 15).Value = 2.2
$ws.Cells.Item(195, 16).Value = 3.2
$ws.Cells.Item(195, 17).Value = 3.5
$ws.Cells.Item(195, 18).Value = -0.25
$ws.Cells.Item(195, 19).Value = 1.85
$ws.Cells.Item(195, 20).Value = 2
$ws.Cells.Item(195, 21).Value = 2.5
$ws.Cells.Item(195, 22).Value = 2.025
$ws.Cells.Item(195, 23).Value = 1.825
$ws.Cells.Item(195, 24).Value = 1.2
$ws.Cells.Item(195, 25).Value = -1
$ws.Cells.Item(195, 26).Value = -1
$ws.Cells.Item(195, 27).Value = 0.8500000000000001
$ws.Cells.Item(195, 28).Value = -1
$ws.Cells.Item(195, 29).Value = -1
$ws.Cells.Item(195, 30).Value = 0.825
